$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACHData")

$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "95125489"
$ws.Range("C5").Value = "95125489"
$ws.Range("D5").Value = "256072691"
$ws.Range("E5").Value = "1"

$ws.Activate()
$ws.Range("E9").Select()
